# The document's images were re-inserted/renamed so that the internal
# drawing-object names no longer match the picture's file extension
# ordering:
#   - header1 (first-page header) BTec logo:     image1.jpg -> image2.jpg
#   - footer  (primary) Pearson logo:             image2.png -> image1.png
#   - footer  (first-page) Pearson logo:          image2.png -> image1.png
#
# InlineShape has no writable .Name in the Word object model, so each
# picture is momentarily converted to a floating Shape (which does expose
# .Name), renamed, then converted back to an inline picture in place.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# wdHeaderFooterIndex constants
$wdHeaderFooterPrimary   = 1
$wdHeaderFooterFirstPage = 2

function Rename-InlineImage($range, [string]$newName) {
    $inline = $range.Range.InlineShapes.Item(1)
    $shape = $inline.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape() | Out-Null
}

# Header (first page) - BTec_Logo-Orange: image1.jpg -> image2.jpg
$hdrFirst = $sec.Headers.Item($wdHeaderFooterFirstPage)
Rename-InlineImage $hdrFirst "image2.jpg"

# Footer (primary/default) - Pearson logo: image2.png -> image1.png
$ftrPrimary = $sec.Footers.Item($wdHeaderFooterPrimary)
Rename-InlineImage $ftrPrimary "image1.png"

# Footer (first page) - Pearson logo: image2.png -> image1.png
$ftrFirst = $sec.Footers.Item($wdHeaderFooterFirstPage)
Rename-InlineImage $ftrFirst "image1.png"
